$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (issue number + reporting week) ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Weekly crime-stat table updates (rows 14-29) ---
$ws.Range("L14").Value = -66.666666666666
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 16
$ws.Range("K15").Value = 14.285714285714
$ws.Range("L15").Value = 14.285714285714
$ws.Range("M15").Value = 220
$ws.Range("N15").Value = -27.272727272727
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 60
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 22.222222222222
$ws.Range("I16").Value = 134
$ws.Range("J16").Value = 111
$ws.Range("K16").Value = 20.72072072072
$ws.Range("L16").Value = 28.846153846153
$ws.Range("M16").Value = 61.44578313253
$ws.Range("N16").Value = -66.749379652605
$ws.Range("C17").Value = 12
$ws.Range("E17").Value = 9.090909090909
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 47
$ws.Range("H17").Value = -8.510638297872
$ws.Range("I17").Value = 199
$ws.Range("J17").Value = 209
$ws.Range("K17").Value = -4.784688995215
$ws.Range("L17").Value = 1.015228426395
$ws.Range("M17").Value = 126.136363636364
$ws.Range("N17").Value = -35.38961038961
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -52.631578947368
$ws.Range("I18").Value = 70
$ws.Range("J18").Value = 112
$ws.Range("K18").Value = -37.5
$ws.Range("L18").Value = 4.477611940298
$ws.Range("M18").Value = 62.790697674418
$ws.Range("N18").Value = -80.978260869565
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 66.666666666666
$ws.Range("I19").Value = 153
$ws.Range("J19").Value = 139
$ws.Range("K19").Value = 10.071942446043
$ws.Range("L19").Value = 25.409836065573
$ws.Range("M19").Value = 104
$ws.Range("N19").Value = 41.666666666666
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 39
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 116.666666666667
$ws.Range("I20").Value = 165
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 65
$ws.Range("L20").Value = 135.714285714286
$ws.Range("M20").Value = 292.857142857143
$ws.Range("N20").Value = -19.117647058823
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = 18.181818181818
$ws.Range("F21").Value = 166
$ws.Range("G21").Value = 134
$ws.Range("H21").Value = 23.880597014925
$ws.Range("I21").Value = 740
$ws.Range("J21").Value = 688
$ws.Range("K21").Value = 7.558139534883
$ws.Range("L21").Value = 26.929674099485
$ws.Range("M21").Value = 117.008797653959
$ws.Range("N21").Value = -47.92399718508
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F22").Value = 1
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 25
$ws.Range("M22").Value = -37.5
$ws.Range("F23").Value = 35
$ws.Range("G23").Value = 36
$ws.Range("H23").Value = -2.777777777777
$ws.Range("I23").Value = 146
$ws.Range("J23").Value = 117
$ws.Range("K23").Value = 24.786324786324
$ws.Range("L23").Value = 124.615384615385
$ws.Range("M23").Value = 124.615384615385
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 76
$ws.Range("G24").Value = 88
$ws.Range("H24").Value = -13.636363636363
$ws.Range("I24").Value = 388
$ws.Range("J24").Value = 377
$ws.Range("K24").Value = 2.917771883289
$ws.Range("L24").Value = 14.117647058823
$ws.Range("M24").Value = 51.5625
$ws.Range("C25").Value = 23
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = 83
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = 20.289855072463
$ws.Range("I25").Value = 401
$ws.Range("J25").Value = 347
$ws.Range("K25").Value = 15.561959654178
$ws.Range("L25").Value = 36.394557823129
$ws.Range("M25").Value = 39.7212543554
$ws.Range("D26").Value = 2
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = -50
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 24
$ws.Range("J26").Value = 23
$ws.Range("K26").Value = 4.347826086956
$ws.Range("L26").Value = 14.285714285714
$ws.Range("C27").Value = 2
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("D27").NumberFormat = "general"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E27").NumberFormat = "general"
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 38
$ws.Range("K27").Value = 58.333333333333
$ws.Range("L27").Value = 90
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -75
$ws.Range("J28").Value = 12
$ws.Range("K28").Value = -16.666666666666
$ws.Range("L28").Value = -47.368421052631
$ws.Range("M28").Value = -37.5
$ws.Range("N28").Value = -64.285714285714
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = -66.666666666666
$ws.Range("J29").Value = 11
$ws.Range("K29").Value = -9.090909090909
$ws.Range("L29").Value = -37.5
$ws.Range("M29").Value = -28.571428571428
$ws.Range("N29").Value = -64.285714285714

